$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("G3").Value = 3.1
$ws.Range("I3").Value = 2.25
$ws.Range("L3").Value = 2.88
$ws.Range("M3").Value = 1.05
$ws.Range("N3").Value = 11
$ws.Range("Q3").Value = 1.93
$ws.Range("R3").Value = 1.88
$ws.Range("W3").Value = 10
$ws.Range("Z3").Value = 34
$ws.Range("AH3").Value = 8
$ws.Range("AM3").Value = 26
$ws.Range("AR3").Value = 81
$ws.Range("AX3").Value = 12
$ws.Range("AY3").Value = 21
$ws.Range("BA3").Value = 51
